$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: add new bug entry (Id=3, bug description, status "À faire")
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Fix l'UI !!! Dans la page translations et subtitles editors lorsqu'on sélectionne un très grand verset"
$ws.Range("C4").Value = "À faire"

# Style C4 like a new "status" fill (theme 5 tint 0.6) with border + center alignment
$ws.Range("C4").Interior.ThemeColor = 6
$ws.Range("C4").Interior.TintAndShade = 0.6
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4108
$ws.Range("C4").Borders.LineStyle = 1

# Row height for row 4
$ws.Rows.Item(4).RowHeight = 48.75

# Clear the style on A22 (remove fill/border)
$ws.Range("A22").Interior.Pattern = -4142
$ws.Range("A22").Borders.LineStyle = -4142
$ws.Range("A22").HorizontalAlignment = -4108
$ws.Range("A22").VerticalAlignment = -4108

# Selection & view changes
$ws.Range("E8").Select()
